# Scheduled market-data refresh: update currentAveragePrice / leve-profit
# columns (H:N) across the per-job sheets with freshly pulled Universalis
# pricing data. Values are plain numbers (no formulas) in the source sheet,
# so we just push the refreshed numbers straight into each cell -- and drop
# the handful of H:N cells whose recompute no longer has data to back them.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2271.4285
$ws.Range("I40").Value = 2300
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 2300
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -2125
$ws.Range("N40").Value = -2550

$ws.Range("H41").Value = 856.05884
$ws.Range("I41").Value = 142.4
$ws.Range("J41").Value = 1153.4166
$ws.Range("K41").Value = 142.4
$ws.Range("L41").Value = 1153.4166
$ws.Range("M41").Value = 297.6
$ws.Range("N41").Value = -2033.4166

$ws.Range("H74").Value = 5360
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 4200
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 4200
$ws.Range("M74").Value = -9064
$ws.Range("N74").Value = -6072

$ws.Range("H77").Value = 5360
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 4200
$ws.Range("K77").Value = 50000
$ws.Range("L77").Value = 21000
$ws.Range("M77").Value = -45320
$ws.Range("N77").Value = -30360

$ws.Range("H86").Value = 1225.75
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 1225.75
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232

$ws.Range("H132").Value = 1222.8649
$ws.Range("I132").Value = 1020.0625
$ws.Range("K132").Value = 3060.1875
$ws.Range("M132").Value = -530.1875

$ws.Range("H137").Value = 23927.455
$ws.Range("I137").Value = 1065.8572
$ws.Range("K137").Value = 3197.5716
$ws.Range("M137").Value = -647.5715999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3654.5312
$ws.Range("I32").Value = 2988.2
$ws.Range("K32").Value = 2988.2
$ws.Range("M32").Value = -2701.2

$ws.Range("H74").Value = 1005.55316
$ws.Range("I74").Value = 845.12195
$ws.Range("K74").Value = 845.12195
$ws.Range("M74").Value = 28.87805000000003

$ws.Range("H77").Value = 1005.55316
$ws.Range("I77").Value = 845.12195
$ws.Range("K77").Value = 4225.60975
$ws.Range("M77").Value = 142.3902500000004

$ws.Range("H97").Value = 775.7692
$ws.Range("I97").Value = 613.1429000000001
$ws.Range("J97").Value = 965.5
$ws.Range("K97").Value = 613.1429000000001
$ws.Range("L97").Value = 965.5
$ws.Range("M97").Value = -117.1429000000001
$ws.Range("N97").Value = -1957.5

$ws.Range("H132").Value = 2394.524
$ws.Range("I132").Value = 1898.8182
$ws.Range("J132").Value = 2939.8
$ws.Range("K132").Value = 5696.4546
$ws.Range("L132").Value = 8819.400000000001
$ws.Range("M132").Value = -3166.4546
$ws.Range("N132").Value = -13879.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 93389.37
$ws.Range("J86").Value = 252887
$ws.Range("L86").Value = 252887
$ws.Range("N86").Value = -255133

$ws.Range("H89").Value = 93389.37
$ws.Range("J89").Value = 252887
$ws.Range("L89").Value = 1264435
$ws.Range("N89").Value = -1275667

$ws.Range("H134").Value = 6070.967
$ws.Range("I134").Value = 6697.269
$ws.Range("K134").Value = 20091.807
$ws.Range("M134").Value = -17556.807

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1773.4166
$ws.Range("I31").Value = 1119.2
$ws.Range("J31").Value = 2240.7144
$ws.Range("K31").Value = 1119.2
$ws.Range("L31").Value = 2240.7144
$ws.Range("M31").Value = -824.2
$ws.Range("N31").Value = -2830.7144

$ws.Range("H34").Value = 1773.4166
$ws.Range("I34").Value = 1119.2
$ws.Range("J34").Value = 2240.7144
$ws.Range("K34").Value = 1119.2
$ws.Range("L34").Value = 2240.7144
$ws.Range("M34").Value = -917.2
$ws.Range("N34").Value = -2644.7144

$ws.Range("H132").Value = 2547.0476
$ws.Range("I132").Value = 1116.091
$ws.Range("K132").Value = 3348.273
$ws.Range("M132").Value = -818.2729999999997

$ws.Range("H134").Value = 1431.6316
$ws.Range("J134").Value = 1000
$ws.Range("L134").Value = 3000
$ws.Range("N134").Value = -8070

$ws.Range("H135").Value = 55617.332
$ws.Range("J135").Value = 55617.332
$ws.Range("L135").Value = 55617.332
$ws.Range("N135").Value = -65757.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1097.5
$ws.Range("I122").Value = 890
$ws.Range("J122").Value = 1149.375
$ws.Range("K122").Value = 8010
$ws.Range("L122").Value = 10344.375
$ws.Range("M122").Value = -5560
$ws.Range("N122").Value = -15244.375

$ws.Range("H129").Value = 24984.467
$ws.Range("J129").Value = 35511.855
$ws.Range("L129").Value = 106535.565
$ws.Range("N129").Value = -116535.565

$ws.Range("H131").Value = 772.53
$ws.Range("J131").Value = 782.84375
$ws.Range("L131").Value = 2348.53125
$ws.Range("N131").Value = -12428.53125

$ws.Range("H140").Value = 1526.75
$ws.Range("I140").Value = 882.8889
$ws.Range("K140").Value = 2648.6667
$ws.Range("M140").Value = 2531.3333

$ws.Range("H141").Value = 2901.2307
$ws.Range("I141").Value = 2901.2307
$ws.Range("K141").Value = 8703.6921
$ws.Range("M141").Value = -3523.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2573480.8
$ws.Range("I126").Value = 2926820
$ws.Range("K126").Value = 8780460
$ws.Range("M126").Value = -8777990

$ws.Range("H132").Value = 3499227
$ws.Range("I132").Value = 5496317
$ws.Range("J132").Value = 4319.25
$ws.Range("K132").Value = 16488951
$ws.Range("L132").Value = 12957.75
$ws.Range("M132").Value = -16486421
$ws.Range("N132").Value = -18017.75

$ws.Range("H138").Value = 49214.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 49214.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 49214.5
$ws.Range("M138").ClearContents() | Out-Null
$ws.Range("N138").Value = -59494.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2034.6666
$ws.Range("I46").Value = 1373
$ws.Range("K46").Value = 1373
$ws.Range("M46").Value = -1185

$ws.Range("H64").Value = 999999
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents() | Out-Null

$ws.Range("H67").Value = 999999
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents() | Out-Null

$ws.Range("H68").Value = 2652.923
$ws.Range("I68").Value = 2478.9
$ws.Range("K68").Value = 2478.9
$ws.Range("M68").Value = -1729.9

$ws.Range("H71").Value = 2652.923
$ws.Range("I71").Value = 2478.9
$ws.Range("K71").Value = 12394.5
$ws.Range("M71").Value = -8650.5

$ws.Range("H122").Value = 11800
$ws.Range("I122").Value = 9666.666999999999
$ws.Range("K122").Value = 29000.001
$ws.Range("M122").Value = -26550.001

$ws.Range("H136").Value = 2829.08
$ws.Range("J136").Value = 5437.375
$ws.Range("L136").Value = 16312.125
$ws.Range("N136").Value = -21412.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 24000
$ws.Range("J63").Value = 24000
$ws.Range("L63").Value = 24000
$ws.Range("N63").Value = -25248

$ws.Range("H66").Value = 24000
$ws.Range("J66").Value = 24000
$ws.Range("L66").Value = 72000
$ws.Range("N66").Value = -78240

$ws.Range("H132").Value = 1447.15
$ws.Range("I132").Value = 1135.7778
$ws.Range("K132").Value = 3407.3334
$ws.Range("M132").Value = -877.3334000000004

$ws.Range("H135").Value = 85888.625
$ws.Range("J135").Value = 85888.625
$ws.Range("L135").Value = 85888.625
$ws.Range("N135").Value = -96028.625
